$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4376.143
$ws.Range("J40").Value = 5635.3335
$ws.Range("L40").Value = 5635.3335
$ws.Range("N40").Value = -5985.3335
$ws.Range("H64").Value = 5319.2607
$ws.Range("I64").Value = 5142.4375
$ws.Range("K64").Value = 5142.4375
$ws.Range("M64").Value = -4894.4375
$ws.Range("H67").Value = 5319.2607
$ws.Range("I67").Value = 5142.4375
$ws.Range("K67").Value = 5142.4375
$ws.Range("M67").Value = -4284.4375
$ws.Range("H74").Value = 5832
$ws.Range("I74").Value = 5998.7334
$ws.Range("J74").Value = 4998.3335
$ws.Range("K74").Value = 5998.7334
$ws.Range("L74").Value = 4998.3335
$ws.Range("M74").Value = -5062.7334
$ws.Range("N74").Value = -6870.3335
$ws.Range("H76").Value = 4679.5
$ws.Range("I76").Value = 3395.3333
$ws.Range("J76").Value = 5730.1816
$ws.Range("K76").Value = 3395.3333
$ws.Range("L76").Value = 5730.1816
$ws.Range("M76").Value = -3080.3333
$ws.Range("N76").Value = -6360.1816
$ws.Range("H77").Value = 5832
$ws.Range("I77").Value = 5998.7334
$ws.Range("J77").Value = 4998.3335
$ws.Range("K77").Value = 29993.667
$ws.Range("L77").Value = 24991.6675
$ws.Range("M77").Value = -25313.667
$ws.Range("N77").Value = -34351.6675
$ws.Range("H79").Value = 4679.5
$ws.Range("I79").Value = 3395.3333
$ws.Range("J79").Value = 5730.1816
$ws.Range("K79").Value = 3395.3333
$ws.Range("L79").Value = 5730.1816
$ws.Range("M79").Value = -2303.3333
$ws.Range("N79").Value = -7914.1816
$ws.Range("H132").Value = 1993.8379
$ws.Range("J132").Value = 4333.3335
$ws.Range("L132").Value = 13000.0005
$ws.Range("N132").Value = -18060.0005
$ws.Range("H137").Value = 2853.1516
$ws.Range("I137").Value = 1674.8
$ws.Range("K137").Value = 5024.4
$ws.Range("M137").Value = -2474.4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 32000
$ws.Range("J39").Value = 32000
$ws.Range("L39").Value = 32000
$ws.Range("N39").Value = -33040
$ws.Range("H63").Value = 2580.1667
$ws.Range("I63").Value = 1257.2727
$ws.Range("J63").Value = 3699.5386
$ws.Range("K63").Value = 1257.2727
$ws.Range("L63").Value = 3699.5386
$ws.Range("M63").Value = -571.2727
$ws.Range("N63").Value = -5071.5386
$ws.Range("H66").Value = 2580.1667
$ws.Range("I66").Value = 1257.2727
$ws.Range("J66").Value = 3699.5386
$ws.Range("K66").Value = 6286.363499999999
$ws.Range("L66").Value = 18497.693
$ws.Range("M66").Value = -2854.363499999999
$ws.Range("N66").Value = -25361.693
$ws.Range("H88").Value = 1529.7222
$ws.Range("I88").Value = 901.6667
$ws.Range("K88").Value = 901.6667
$ws.Range("M88").Value = -495.6667
$ws.Range("H91").Value = 1529.7222
$ws.Range("I91").Value = 901.6667
$ws.Range("K91").Value = 901.6667
$ws.Range("M91").Value = 502.3333
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1214.3077
$ws.Range("J20").Value = 1524.6666
$ws.Range("L20").Value = 1524.6666
$ws.Range("N20").Value = -2018.6666
$ws.Range("H54").Value = 28599.2
$ws.Range("I54").Value = 7999
$ws.Range("J54").Value = 33749.25
$ws.Range("K54").Value = 7999
$ws.Range("L54").Value = 33749.25
$ws.Range("M54").Value = -7515
$ws.Range("N54").Value = -34717.25
$ws.Range("H99").Value = 2415.52
$ws.Range("I99").Value = 2161.476
$ws.Range("K99").Value = 2161.476
$ws.Range("M99").Value = -663.4760000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2800.8333
$ws.Range("J16").Value = 2361
$ws.Range("L16").Value = 2361
$ws.Range("N16").Value = -2935
$ws.Range("H58").Value = 3692.4614
$ws.Range("I58").Value = 1262.75
$ws.Range("K58").Value = 1262.75
$ws.Range("M58").Value = -1059.75
$ws.Range("H99").Value = 2965.5476
$ws.Range("I99").Value = 2914.4688
$ws.Range("J99").Value = 3129
$ws.Range("K99").Value = 2914.4688
$ws.Range("L99").Value = 3129
$ws.Range("M99").Value = -1416.4688
$ws.Range("N99").Value = -6125
$ws.Range("H113").Value = 2800.8333
$ws.Range("J113").Value = 2361
$ws.Range("L113").Value = 2361
$ws.Range("N113").Value = -6701
$ws.Range("H126").Value = 2965.5476
$ws.Range("I126").Value = 2914.4688
$ws.Range("J126").Value = 3129
$ws.Range("K126").Value = 8743.4064
$ws.Range("L126").Value = 9387
$ws.Range("M126").Value = -6273.4064
$ws.Range("N126").Value = -14327
$ws.Range("H136").Value = 3692.4614
$ws.Range("I136").Value = 1262.75
$ws.Range("K136").Value = 3788.25
$ws.Range("M136").Value = -1238.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 556.625
$ws.Range("I8").Value = 556.625
$ws.Range("K8").Value = 1669.875
$ws.Range("M8").Value = -1530.875
$ws.Range("H23").Value = 105.71429
$ws.Range("I23").Value = 63.5
$ws.Range("J23").Value = 162
$ws.Range("K23").Value = 190.5
$ws.Range("L23").Value = 486
$ws.Range("M23").Value = 44.5
$ws.Range("N23").Value = -956
$ws.Range("H75").Value = 4467.778
$ws.Range("I75").Value = 3850
$ws.Range("J75").Value = 4644.2856
$ws.Range("K75").Value = 11550
$ws.Range("L75").Value = 13932.8568
$ws.Range("M75").Value = -10552
$ws.Range("N75").Value = -15928.8568
$ws.Range("H78").Value = 4467.778
$ws.Range("I78").Value = 3850
$ws.Range("J78").Value = 4644.2856
$ws.Range("K78").Value = 34650
$ws.Range("L78").Value = 41798.5704
$ws.Range("M78").Value = -29658
$ws.Range("N78").Value = -51782.5704
$ws.Range("H92").Value = 1499
$ws.Range("I92").Value = 1499
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 4497
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -3249
$ws.Range("N92").ClearContents()
$ws.Range("H116").Value = 20000
$ws.Range("J116").Value = 20000
$ws.Range("L116").Value = 60000
$ws.Range("N116").Value = -66884
$ws.Range("H141").Value = 3139.1482
$ws.Range("I141").Value = 2384.1538
$ws.Range("J141").Value = 3840.2144
$ws.Range("K141").Value = 7152.4614
$ws.Range("L141").Value = 11520.6432
$ws.Range("M141").Value = -1972.4614
$ws.Range("N141").Value = -21880.6432
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 756.1852
$ws.Range("I97").Value = 316.05264
$ws.Range("J97").Value = 995.1142599999999
$ws.Range("K97").Value = 316.05264
$ws.Range("L97").Value = 995.1142599999999
$ws.Range("M97").Value = 179.94736
$ws.Range("N97").Value = -1987.11426
$ws.Range("H113").Value = 3432.6365
$ws.Range("I113").Value = 2866.5625
$ws.Range("J113").Value = 4942.1665
$ws.Range("K113").Value = 2866.5625
$ws.Range("L113").Value = 4942.1665
$ws.Range("M113").Value = -696.5625
$ws.Range("N113").Value = -9282.166499999999
$ws.Range("H122").Value = 4460.0835
$ws.Range("I122").Value = 2395.25
$ws.Range("K122").Value = 7185.75
$ws.Range("M122").Value = -4735.75
$ws.Range("H132").Value = 4100.107
$ws.Range("I132").Value = 2660.2354
$ws.Range("K132").Value = 7980.706200000001
$ws.Range("M132").Value = -5450.706200000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3448.647
$ws.Range("I61").Value = 1610.2307
$ws.Range("K61").Value = 1610.2307
$ws.Range("M61").Value = -1408.2307
$ws.Range("H74").Value = 128129.3
$ws.Range("I74").Value = 54938
$ws.Range("J74").Value = 146427.12
$ws.Range("K74").Value = 54938
$ws.Range("L74").Value = 146427.12
$ws.Range("M74").Value = -53940
$ws.Range("N74").Value = -148423.12
$ws.Range("H77").Value = 128129.3
$ws.Range("I77").Value = 54938
$ws.Range("J77").Value = 146427.12
$ws.Range("K77").Value = 164814
$ws.Range("L77").Value = 439281.36
$ws.Range("M77").Value = -159822
$ws.Range("N77").Value = -449265.36
$ws.Range("H80").Value = 96500
$ws.Range("J80").Value = 96500
$ws.Range("L80").Value = 96500
$ws.Range("N80").Value = -98746
$ws.Range("H83").Value = 96500
$ws.Range("J83").Value = 96500
$ws.Range("L83").Value = 289500
$ws.Range("N83").Value = -300732
$ws.Range("H113").Value = 3448.647
$ws.Range("I113").Value = 1610.2307
$ws.Range("K113").Value = 1610.2307
$ws.Range("M113").Value = 559.7692999999999
$ws.Range("H122").Value = 4917.5
$ws.Range("I122").Value = 4207.4287
$ws.Range("K122").Value = 12622.2861
$ws.Range("M122").Value = -10172.2861
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 8472
$ws.Range("I52").Value = 4599
$ws.Range("K52").Value = 4599
$ws.Range("M52").Value = -4373
$ws.Range("H74").Value = 18683.334
$ws.Range("J74").Value = 18683.334
$ws.Range("L74").Value = 18683.334
$ws.Range("N74").Value = -20555.334
$ws.Range("H77").Value = 18683.334
$ws.Range("J77").Value = 18683.334
$ws.Range("L77").Value = 56050.00199999999
$ws.Range("N77").Value = -65410.00199999999
$ws.Range("H107").Value = 1001176.5
$ws.Range("I107").Value = 1819417.1
$ws.Range("J107").Value = 1104.5555
$ws.Range("K107").Value = 5458251.300000001
$ws.Range("L107").Value = 3313.6665
$ws.Range("M107").Value = -5456331.300000001
$ws.Range("N107").Value = -7153.666499999999
$ws.Range("H109").Value = 143272.12
$ws.Range("J109").Value = 143272.12
$ws.Range("L109").Value = 143272.12
$ws.Range("N109").Value = -146046.12
$ws.Range("H113").Value = 892.4783
$ws.Range("J113").Value = 1316.2
$ws.Range("L113").Value = 3948.6
$ws.Range("N113").Value = -8288.6
$ws.Range("H132").Value = 3038.182
$ws.Range("I132").Value = 2260.5264
$ws.Range("K132").Value = 6781.5792
$ws.Range("M132").Value = -4251.5792
